$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# --- Row 22 (fm17): populate the previously-blank test row with the new
#     "% Loss deductible with min and max deductible" (1st level) test ---

# Copy formatting from analogous rows first so number formats / alignment
# match the rest of the table (D-column needs the right-aligned "All" style,
# E-column needs the plain/default style used by similar numeric Calcrule cells).
$ws.Range("D8").Copy()
$ws.Range("D22").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E40").Copy()
$ws.Range("E22").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C22").Value2 = "% Loss deductible with min and max deductible"
$ws.Range("D22").Value2 = "All"
$ws.Range("E22").Value2 = 19
$ws.Range("F22").Value2 = 1
$ws.Range("G22").Value2 = 1
$ws.Range("H22").Value2 = "complete"
$ws.Range("I22").Value2 = "complete"

# --- Row 40 (fm35): turn the former single-level test into the "2nd level"
#     test, now applicable to all allocrules and spanning 2 levels ---
$ws.Range("C40").Value2 = "% Loss deductible with min and max deductible 2nd level. Calcrule 19"
$ws.Range("D40").Value2 = "All"
$ws.Range("F40").Value2 = 2

# --- Leave the selection where the author left it when they saved ---
$ws.Activate() | Out-Null
$ws.Range("D41").Select() | Out-Null
